$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New daily data rows appended at the bottom of the table (rows 191-193),
# matching the "Actualización desde MV -datos-" update.
$newRows = @(
    @{ Row = 191; Date = "04-10-2021"; B = 3166; C = 7410 },
    @{ Row = 192; Date = "05-10-2021"; B = 3584; C = 7875 },
    @{ Row = 193; Date = "06-10-2021"; B = 4005; C = 9176 }
)

foreach ($r in $newRows) {
    $rowIndex = $r.Row

    # Force the date column to be stored as plain text (matching the rest
    # of column A), then drop the formatting again so the cell keeps the
    # default (unstyled) look used by all the other data rows.
    $cellA = $ws.Cells.Item($rowIndex, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $r.Date
    $cellA.ClearFormats()

    $ws.Cells.Item($rowIndex, 2).Value = $r.B
    $ws.Cells.Item($rowIndex, 3).Value = $r.C
}
